# Apply updated cryptocurrency price / volume(1h) data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.699.52"
$ws.Range("E2").Value = "  +3.31%  "
$ws.Range("D3").Value = "3.293.93"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  +3.22%  "
$ws.Range("D9").Value = "3.285.45"
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("E10").Value = "  +1.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.576"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.79"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("E13").Value = "  +3.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "695.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +14.24%  "
$ws.Range("D15").Value = "3.821.50"
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("E16").Value = "  +0.71%  "
$ws.Range("D17").Value = "67.828.55"
$ws.Range("E17").Value = "  +3.33%  "
$ws.Range("E18").Value = "  +1.51%  "
$ws.Range("D19").Value = "3.292.94"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.69%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.897"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("E27").Value = "  +1.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.78%  "
$ws.Range("E30").Value = "  +2.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "584.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.28%  "
$ws.Range("D34").Value = "3.869.75"
$ws.Range("E34").Value = "  +1.99%  "
$ws.Range("E35").Value = "  +1.94%  "
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("E37").Value = "  -8.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.00%  "
$ws.Range("E39").Value = "  +1.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.63%  "
$ws.Range("E41").Value = "  +2.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "32.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").Value = "0.0₃0682"
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.332"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0413"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.49%  "
$ws.Range("E47").Value = "  +2.45%  "
$ws.Range("E48").Value = "  +10.48%  "
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "129.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.02%  "
